# Regenerate merged AHB files
# - Rename header labels from _old/_new suffix to _FV2304/_FV2310 suffix
# - Add a table (ListObject) over the data range A1:U60
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldCols = @('Segmentname','Segmentgruppe','Segment','Datenelement','Segment ID','Code','Qualifier','Beschreibung','Bedingungsausdruck','Bedingung')

# Columns A-J: "<name>_old" -> "<name>_FV2304"
for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldCols[$i] + "_FV2304"
}

# Column K stays "diff" (unchanged)

# Columns L-U: "<name>_new" -> "<name>_FV2310"
for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $oldCols[$i] + "_FV2310"
}

# Add a table over the full used range
$range = $ws.Range("A1:U60")
$table = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$table.Name = "Table1"

# Freeze the header row
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
